# Generate Report for Handoff
#
# This updates the localization-status report:
#  - The "Priority" column (E) for the last four rows of both the
#    zh-cn and de-de sheets changes from "low" to "ht" (content is no
#    longer a duplicate, so the priority is recalculated).
#  - The zh-cn sheet's "Latest Handoff Datetime" (H) for the same four
#    rows is refreshed to a later timestamp.
#  - The "Ready for handoff" timestamp (Overview!G4:G7 and
#    de-de!H4:H7, which share one string) is refreshed to a later
#    timestamp as well.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Priority column: "low" -> "ht" for rows 4-7 on both locale sheets.
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("E5").Value = "ht"
$wsZhCn.Range("E6").Value = "ht"
$wsZhCn.Range("E7").Value = "ht"

$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("E5").Value = "ht"
$wsDeDe.Range("E6").Value = "ht"
$wsDeDe.Range("E7").Value = "ht"

# zh-cn Latest Handoff Datetime refresh for rows 4-7.
$wsZhCn.Range("H4").Value = "2016-08-31 10:34:57"
$wsZhCn.Range("H5").Value = "2016-08-31 10:34:57"
$wsZhCn.Range("H6").Value = "2016-08-31 10:34:57"
$wsZhCn.Range("H7").Value = "2016-08-31 10:34:57"

# "Ready for handoff" timestamp refresh (Overview G4:G7, de-de H4:H7).
$wsOverview.Range("G4").Value = "2016-08-31 10:35:11"
$wsOverview.Range("G5").Value = "2016-08-31 10:35:11"
$wsOverview.Range("G6").Value = "2016-08-31 10:35:11"
$wsOverview.Range("G7").Value = "2016-08-31 10:35:11"

$wsDeDe.Range("H4").Value = "2016-08-31 10:35:11"
$wsDeDe.Range("H5").Value = "2016-08-31 10:35:11"
$wsDeDe.Range("H6").Value = "2016-08-31 10:35:11"
$wsDeDe.Range("H7").Value = "2016-08-31 10:35:11"
